# Applies the cryptos.xlsx price/volume refresh described in the commit
# "Updated cryptos list on Wed Apr 26 07:59:52 UTC 2023 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.395.91"
$ws.Range("E2").Value = "  +3.51%  "

$ws.Range("D3").Value = "1.864.00"
$ws.Range("E3").Value = "  +2.03%  "

$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  -0.49%  "

$ws.Range("D5").Value = "'337.05"
$ws.Range("E5").Value = "  +1.75%  "

$ws.Range("E6").Value = "  -0.54%  "

$ws.Range("D7").Value = "'0.4688"
$ws.Range("E7").Value = "  +2.34%  "

$ws.Range("D8").Value = "'0.3970"
$ws.Range("E8").Value = "  +3.94%  "

$ws.Range("D9").Value = "'47.56"
$ws.Range("E9").Value = "  +2.70%  "

$ws.Range("D10").Value = "'0.08020"
$ws.Range("E10").Value = "  +1.51%  "

$ws.Range("D11").Value = "'0.9943"
$ws.Range("E11").Value = "  +2.58%  "

$ws.Range("D12").Value = "'21.93"
$ws.Range("E12").Value = "  +4.26%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.871.09"
$ws.Range("E13").Value = "  +1.88%  "

$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "'6.022"
$ws.Range("E14").Value = "  +2.62%  "

$ws.Range("D15").Value = "'7.249"
$ws.Range("E15").Value = "  +3.06%  "

$ws.Range("D16").Value = "'90.21"
$ws.Range("E16").Value = "  +2.34%  "

$ws.Range("E17").Value = "  -0.47%  "

$ws.Range("D18").Value = "'0.00001039"
$ws.Range("E18").Value = "  +0.90%  "

$ws.Range("D19").Value = "'0.06603"
$ws.Range("E19").Value = "  -0.73%  "

$ws.Range("D20").Value = "'17.50"
$ws.Range("E20").Value = "  +2.08%  "

$ws.Range("E21").Value = "  -0.45%  "

$ws.Range("D22").Value = "28.405.79"
$ws.Range("E22").Value = "  +3.55%  "

$ws.Range("D23").Value = "'5.464"
$ws.Range("E23").Value = "  +2.83%  "

$ws.Range("D24").Value = "'11.01"

$ws.Range("D25").Value = "'2.268"
$ws.Range("E25").Value = "  -1.81%  "

$ws.Range("D26").Value = "2.085.78"
$ws.Range("E26").Value = "  +1.56%  "

$ws.Range("D27").Value = "'161.11"
$ws.Range("E27").Value = "  +2.56%  "

$ws.Range("D28").Value = "'19.71"
$ws.Range("E28").Value = "  +1.85%  "

$ws.Range("D29").Value = "'2.107"
$ws.Range("E29").Value = "  +2.32%  "

$ws.Range("D30").Value = "'5.455"
$ws.Range("E30").Value = "  +4.36%  "

$ws.Range("D31").Value = "'119.19"
$ws.Range("E31").Value = "  +0.85%  "

$ws.Range("D32").Value = "'0.09507"
$ws.Range("E32").Value = "  +2.27%  "

$ws.Range("D33").Value = "'0.9599"
$ws.Range("E33").Value = "  +1.66%  "

$ws.Range("D34").Value = "'3.592"

$ws.Range("D35").Value = "'5.348"
$ws.Range("E35").Value = "  +2.21%  "

$ws.Range("D36").Value = "'1.373"
$ws.Range("E36").Value = "  +4.76%  "

$ws.Range("D37").Value = "'0.06154"
$ws.Range("E37").Value = "  +3.79%  "

$ws.Range("D38").Value = "'0.02244"
$ws.Range("E38").Value = "  +2.75%  "

$ws.Range("D39").Value = "'8.282"
$ws.Range("E39").Value = "  +3.92%  "

$ws.Range("D40").Value = "'1.177"
$ws.Range("E40").Value = "  +1.85%  "

$ws.Range("D41").Value = "'0.5912"
$ws.Range("E41").Value = "  +2.29%  "

$ws.Range("E42").Value = "  -0.51%  "

$ws.Range("D43").Value = "'0.1871"
$ws.Range("E43").Value = "  +2.00%  "

$ws.Range("D44").Value = "'10.28"
$ws.Range("E44").Value = "  +2.78%  "

$ws.Range("D45").Value = "'1.274"
$ws.Range("E45").Value = "  -0.02%  "

$ws.Range("D46").Value = "'0.07592"
$ws.Range("E46").Value = "  +14.39%  "

$ws.Range("D47").Value = "'0.5534"

$ws.Range("D48").Value = "'12.09"
$ws.Range("E48").Value = "  +1.18%  "

$ws.Range("D49").Value = "'1.940"
$ws.Range("E49").Value = "  +4.06%  "

$ws.Range("D50").Value = "'2.064"
$ws.Range("E50").Value = "  +13.25%  "

$ws.Range("D51").Value = "'111.89"
$ws.Range("E51").Value = "  +1.89%  "
